$d = $word.ActiveDocument

# --- Hunk 1: remove the old "_GoBack" bookmark (was after "...efektif") ---
$d.Bookmarks("_GoBack").Delete()

# --- Hunk 2: " data " -> " prangko dan filateli, meliputi data " plus a
#     relocated "_GoBack" bookmark right after "filateli," ---
$rng = $d.Content
$rng.Find.Execute("adalah data transaksi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorStart = $rng.Start

# "adalah" is 6 characters; " data " (the text being replaced) starts right after it.
$dataStart = $anchorStart + 6
$dataEnd = $dataStart + 6

$target = $d.Range($dataStart, $dataEnd)
$target.Text = ""

# Build up the replacement text piece by piece so each piece keeps its own
# run (matching how the real edit was typed/spell-checked word by word),
# using short-lived "wedge" bookmarks to keep the runs from coalescing.
$pieces = @(" ", "prangko", " ", "dan", " ", "filateli", ",")
$pos = $dataStart
$wedgeCount = 0
for ($i = 0; $i -lt $pieces.Length; $i++) {
    $piece = $pieces[$i]
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter($piece)
    $pos = $pos + $piece.Length
    $wedge = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_wedge" + $wedgeCount, $wedge)
    $wedgeCount = $wedgeCount + 1
}

# The real "_GoBack" bookmark belongs exactly here, right after "filateli,".
$d.Bookmarks.Add("_GoBack", $d.Range($pos, $pos))

$pieces2 = @(" ", "meliputi", " ", "data ")
for ($i = 0; $i -lt $pieces2.Length; $i++) {
    $piece = $pieces2[$i]
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter($piece)
    $pos = $pos + $piece.Length
    if ($i -lt $pieces2.Length - 1) {
        $wedge = $d.Range($pos, $pos)
        $d.Bookmarks.Add("_wedge" + $wedgeCount, $wedge)
        $wedgeCount = $wedgeCount + 1
    }
}

# Drop the temporary wedges now that every piece has been typed - the runs
# they separated stay split even once the bookmark itself is gone.
for ($i = 0; $i -lt $wedgeCount; $i++) {
    $d.Bookmarks("_wedge" + $i).Delete()
}
